$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the data table by two rows (95 -> 97), mirroring the existing
# "NTC" control-row pattern used elsewhere in the sheet. Copy the format +
# values of row 94 (a blank/0 sample row) down into the two new rows, then
# fix up the running index (col A) and sample name (col B -> "NTC").
$ws.Range("A94:H94").Copy($ws.Range("A96:H96"))
$ws.Range("A94:H94").Copy($ws.Range("A97:H97"))

$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = "NTC"

$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = "NTC"

# --- Update the print area to include the new rows, pushing the previous
# print-area history down one slot (mirrors the workbook's existing
# Print_Area / Print_Area_0 / Print_Area_0_0... history chain).
$ws.Names.Item("_xlnm.Print_Area_0").RefersTo = '=_!$A$1:$H$95'
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0", '=_!$A$1:$H$87')
$ws.PageSetup.PrintArea = '$A$1:$H$97'

# --- Shrink the print scale slightly to keep everything fitting one page.
$ws.PageSetup.Zoom = 52

# --- Move the selection to reflect where editing left off.
$null = $ws.Range("B98").Select()
